$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.920.61"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.213.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.34%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.212.98"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.32%  "

$ws.Range("E9").Value = "  +0.21%  "

$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.508"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.741.27"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.36%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.111.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.70%  "

$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.212.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.31%  "

$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.32"
$ws.Range("D20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.83%  "

$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  +2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("E29").Value = "  +3.45%  "

$ws.Range("E30").Value = "  +3.32%  "

$ws.Range("E31").Value = "  +7.58%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "

$ws.Range("E33").Value = "  +1.47%  "

$ws.Range("E34").Value = "  +0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.96%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0903"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "477.27"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.90%  "

$ws.Range("E39").Value = "  +0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.19%  "

$ws.Range("E42").Value = "  +4.19%  "

$ws.Range("E43").Value = "  +0.50%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.953.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.45"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0641"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.22%  "

$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("E48").Value = "  +0.09%  "

$ws.Range("E49").Value = "  +0.82%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.50%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.79%  "
